# edit.ps1
# Applies the cryptos list price/volume refresh described by the commit:
# "Updated cryptos list on Sat Sep 14 19:56:57 UTC 2024 with GitHub Actions"
#
# Most cells are plain display strings (e.g. "  +0.05%  ", "59.832.80") and
# can be written directly via Range.Value. A handful of Price (column D)
# values look like plain numbers (e.g. "551.23", "1.00"); Excel's normal
# Range.Value setter auto-converts those to numeric cells, which would lose
# the original text formatting. Set-TextValue forces those through as text
# (quote-prefix trick) and then restores the cell's style to "Normal" so no
# stray number-format/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# Plain string/text assignments (not numeric-looking, safe to set directly)
$ws.Range("D2").Value = '59.832.80'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.409.37'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  -2.44%  '
$ws.Range("E11").Value = '  -0.97%  '
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("E13").Value = '  +2.53%  '
$ws.Range("D14").Value = '2.838.88'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = '59.766.35'
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").Value = '2.424.29'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("E21").Value = '  -3.79%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("E23").Value = '  +2.46%  '
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").Value = '0.0₃0770'
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("E29").Value = '  -1.89%  '
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  -4.16%  '
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("E35").Value = '  -0.45%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  -1.63%  '
$ws.Range("E38").Value = '  -1.83%  '
$ws.Range("E39").Value = '  +2.27%  '
$ws.Range("E40").Value = '  -4.29%  '
$ws.Range("E41").Value = '  -2.10%  '
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("E44").Value = '  +1.29%  '
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("B47").Value = 'Polygon'
$ws.Range("C47").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("E47").Value = '  -4.32%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  -3.03%  '
$ws.Range("E51").Value = '  -1.23%  '

# Numeric-looking Price values: force as text via quote-prefix technique,
# then reset style back to Normal so no residual formatting is left on the cell.
Set-TextValue $ws.Range("D5") '551.23'
Set-TextValue $ws.Range("D6") '136.86'
Set-TextValue $ws.Range("D12") '0.353'
Set-TextValue $ws.Range("D13") '25.29'
Set-TextValue $ws.Range("D20") '328.70'
Set-TextValue $ws.Range("D22") '1.00'
Set-TextValue $ws.Range("D23") '66.22'
Set-TextValue $ws.Range("D30") '169.06'
Set-TextValue $ws.Range("D39") '320.47'
Set-TextValue $ws.Range("D40") '0.404'
Set-TextValue $ws.Range("D42") '139.89'
Set-TextValue $ws.Range("D44") '19.46'
Set-TextValue $ws.Range("D47") '0.397'
Set-TextValue $ws.Range("D48") '0.0222'
Set-TextValue $ws.Range("D49") '11.04'
